$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"3.820521676800581e-11"
$ws.Range("C2").Value = [double]"2.22304730179701e-09"
$ws.Range("D2").Value = [double]"337.1190423067083"
$ws.Range("E2").Value = [double]"2367095152636972"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = [double]"2367095152637308"
